$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextCell $ws.Range("D2") "30.314.56"
$ws.Range("E2").Value = "  +0.03%  "

Set-TextCell $ws.Range("D3") "1.929.88"
$ws.Range("E3").Value = "  +0.01%  "

$ws.Range("E4").Value = "  +0.37%  "

Set-TextCell $ws.Range("D5") "0.7433"
$ws.Range("E5").Value = "  +3.06%  "

$ws.Range("E6").Value = "  -2.25%  "

$ws.Range("E7").Value = "  +0.35%  "

$ws.Range("B8").Value = "Cardano"
$ws.Range("C8").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
Set-TextCell $ws.Range("D8") "0.3154"
$ws.Range("E8").Value = "  -1.91%  "

$ws.Range("B9").Value = "Solana"
$ws.Range("C9").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
Set-TextCell $ws.Range("D9") "27.43"
$ws.Range("E9").Value = "  -2.02%  "

Set-TextCell $ws.Range("D10") "0.07035"
$ws.Range("E10").Value = "  -1.00%  "

Set-TextCell $ws.Range("D11") "0.08032"
$ws.Range("E11").Value = "  +0.30%  "

Set-TextCell $ws.Range("D12") "0.7768"

Set-TextCell $ws.Range("D13") "1.912.44"
$ws.Range("E13").Value = "  -0.92%  "

$ws.Range("E14").Value = "  -0.15%  "

Set-TextCell $ws.Range("D15") "93.13"
$ws.Range("E15").Value = "  -1.79%  "

$ws.Range("E16").Value = "  -1.34%  "

Set-TextCell $ws.Range("D17") "30.317.19"
$ws.Range("E17").Value = "  +0.07%  "

Set-TextCell $ws.Range("D18") "5.981"
$ws.Range("E18").Value = "  +3.93%  "

Set-TextCell $ws.Range("D19") "250.47"
$ws.Range("E19").Value = "  -2.57%  "

Set-TextCell $ws.Range("D20") "0.000007926"
$ws.Range("E20").Value = "  -2.06%  "

$ws.Range("E21").Value = "  +0.27%  "

Set-TextCell $ws.Range("D22") "2.149.95"
$ws.Range("E22").Value = "  -1.43%  "

Set-TextCell $ws.Range("D23") "1.004"
$ws.Range("E23").Value = "  +0.29%  "

Set-TextCell $ws.Range("D24") "6.647"
$ws.Range("E24").Value = "  -2.71%  "

Set-TextCell $ws.Range("D25") "9.558"
$ws.Range("E25").Value = "  -0.19%  "

Set-TextCell $ws.Range("D26") "165.50"
$ws.Range("E26").Value = "  +0.48%  "

Set-TextCell $ws.Range("D27") "19.02"
$ws.Range("E27").Value = "  -0.60%  "

$ws.Range("E28").Value = "  -0.14%  "

Set-TextCell $ws.Range("D29") "2.169"

Set-TextCell $ws.Range("D30") "1.562"
$ws.Range("E30").Value = "  +1.94%  "

$ws.Range("E31").Value = "  -0.13%  "

Set-TextCell $ws.Range("D32") "4.414"
$ws.Range("E32").Value = "  -0.13%  "

Set-TextCell $ws.Range("D33") "4.108"
$ws.Range("E33").Value = "  -1.10%  "

Set-TextCell $ws.Range("D34") "0.05217"
$ws.Range("E34").Value = "  +1.71%  "

Set-TextCell $ws.Range("D35") "1.311"
$ws.Range("E35").Value = "  +2.04%  "

Set-TextCell $ws.Range("D36") "0.7541"
$ws.Range("E36").Value = "  +0.57%  "

Set-TextCell $ws.Range("D37") "2.767"
$ws.Range("E37").Value = "  -0.16%  "

Set-TextCell $ws.Range("D38") "0.01952"
$ws.Range("E38").Value = "  -1.78%  "

Set-TextCell $ws.Range("D39") "2.789"
$ws.Range("E39").Value = "  -0.31%  "

Set-TextCell $ws.Range("D40") "6.523"
$ws.Range("E40").Value = "  +1.98%  "

Set-TextCell $ws.Range("D41") "76.61"
$ws.Range("E41").Value = "  -2.25%  "

$ws.Range("E42").Value = "  -0.65%  "

Set-TextCell $ws.Range("D43") "1.954"
$ws.Range("E43").Value = "  -2.35%  "

Set-TextCell $ws.Range("D44") "0.8426"
$ws.Range("E44").Value = "  -0.39%  "

Set-TextCell $ws.Range("D45") "1.002"
$ws.Range("E45").Value = "  +0.32%  "

Set-TextCell $ws.Range("D46") "7.672"

Set-TextCell $ws.Range("D47") "9.954"
$ws.Range("E47").Value = "  +1.32%  "

Set-TextCell $ws.Range("D48") "101.15"
$ws.Range("E48").Value = "  +0.02%  "

$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
Set-TextCell $ws.Range("D49") "37.55"
$ws.Range("E49").Value = "  +1.74%  "

$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
Set-TextCell $ws.Range("D50") "2.060.38"
$ws.Range("E50").Value = "  -1.66%  "

Set-TextCell $ws.Range("D51") "0.1227"
$ws.Range("E51").Value = "  +7.20%  "
